# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Wed Jun 19 20:41:32 UTC 2024 with GitHub Actions".
# D-column values that are plain numbers are written with a leading
# apostrophe so Excel stores them as text (matching the sheet's existing
# text-typed Price column) instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '64.829.50'
$ws.Cells.Item(2, 5).Value = '  +0.69%  '

$ws.Cells.Item(3, 4).Value = '3.551.56'
$ws.Cells.Item(3, 5).Value = '  +3.86%  '

$ws.Cells.Item(4, 4).Value = '''1.00'
$ws.Cells.Item(4, 5).Value = '  +0.00%  '

$ws.Cells.Item(5, 4).Value = '''600.28'
$ws.Cells.Item(5, 5).Value = '  +3.34%  '

$ws.Cells.Item(6, 4).Value = '''135.60'
$ws.Cells.Item(6, 5).Value = '  +1.02%  '

$ws.Cells.Item(7, 4).Value = '3.550.78'
$ws.Cells.Item(7, 5).Value = '  +3.87%  '

$ws.Cells.Item(8, 5).Value = '  -0.02%  '

$ws.Cells.Item(9, 4).Value = '''0.495'
$ws.Cells.Item(9, 5).Value = '  +2.64%  '

$ws.Cells.Item(10, 5).Value = '  +1.83%  '

$ws.Cells.Item(11, 4).Value = '''6.91'
$ws.Cells.Item(11, 5).Value = '  -1.09%  '

$ws.Cells.Item(12, 4).Value = '''0.386'
$ws.Cells.Item(12, 5).Value = '  +3.13%  '

$ws.Cells.Item(13, 4).Value = '4.156.17'
$ws.Cells.Item(13, 5).Value = '  +3.89%  '

$ws.Cells.Item(14, 5).Value = '  +2.22%  '

$ws.Cells.Item(15, 4).Value = '3.557.05'
$ws.Cells.Item(15, 5).Value = '  +3.70%  '

$ws.Cells.Item(16, 4).Value = '''26.97'
$ws.Cells.Item(16, 5).Value = '  +3.15%  '

$ws.Cells.Item(18, 4).Value = '64.727.12'
$ws.Cells.Item(18, 5).Value = '  +0.49%  '

$ws.Cells.Item(19, 4).Value = '''10.02'
$ws.Cells.Item(19, 5).Value = '  +4.62%  '

$ws.Cells.Item(20, 5).Value = '  +6.25%  '

$ws.Cells.Item(21, 5).Value = '  +3.04%  '

$ws.Cells.Item(22, 4).Value = '''386.79'
$ws.Cells.Item(22, 5).Value = '  +2.14%  '

$ws.Cells.Item(23, 4).Value = '''0.576'
$ws.Cells.Item(23, 5).Value = '  +6.29%  '

$ws.Cells.Item(24, 4).Value = '3.694.44'
$ws.Cells.Item(24, 5).Value = '  +3.75%  '

$ws.Cells.Item(25, 4).Value = '''74.36'
$ws.Cells.Item(25, 5).Value = '  +3.54%  '

$ws.Cells.Item(26, 5).Value = '  +0.13%  '

$ws.Cells.Item(27, 4).Value = '''0.0000117'
$ws.Cells.Item(27, 5).Value = '  +11.93%  '

$ws.Cells.Item(28, 4).Value = '''7.63'
$ws.Cells.Item(28, 5).Value = '  +7.43%  '

$ws.Cells.Item(29, 4).Value = '''1.00'
$ws.Cells.Item(29, 5).Value = '  +0.02%  '

$ws.Cells.Item(30, 4).Value = '''2.30'
$ws.Cells.Item(30, 5).Value = '  +5.39%  '

$ws.Cells.Item(31, 5).Value = '  +4.10%  '

$ws.Cells.Item(32, 4).Value = '3.559.92'
$ws.Cells.Item(32, 5).Value = '  +3.65%  '

$ws.Cells.Item(33, 5).Value = '  +23.48%  '

$ws.Cells.Item(34, 2).Value = 'USDe'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(34, 4).Value = '''1.00'
$ws.Cells.Item(34, 5).Value = '  +0.02%  '

$ws.Cells.Item(35, 2).Value = 'EthereumClassic'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(35, 4).Value = '''23.92'
$ws.Cells.Item(35, 5).Value = '  +4.37%  '

$ws.Cells.Item(36, 5).Value = '  +2.56%  '

$ws.Cells.Item(37, 4).Value = '''170.00'
$ws.Cells.Item(37, 5).Value = '  -0.57%  '

$ws.Cells.Item(38, 4).Value = '''6.92'

$ws.Cells.Item(39, 5).Value = '  +6.07%  '

$ws.Cells.Item(40, 4).Value = '''4.99'
$ws.Cells.Item(40, 5).Value = '  +7.68%  '

$ws.Cells.Item(41, 4).Value = '''0.0805'
$ws.Cells.Item(41, 5).Value = '  +6.27%  '

$ws.Cells.Item(42, 4).Value = '''0.827'
$ws.Cells.Item(42, 5).Value = '  +3.50%  '

$ws.Cells.Item(43, 4).Value = '''26.93'
$ws.Cells.Item(43, 5).Value = '  +20.10%  '

$ws.Cells.Item(44, 5).Value = '  +2.44%  '

$ws.Cells.Item(45, 4).Value = '''1.00'
$ws.Cells.Item(45, 5).Value = '  -0.05%  '

$ws.Cells.Item(46, 5).Value = '  +4.61%  '

$ws.Cells.Item(47, 5).Value = '  +10.27%  '

$ws.Cells.Item(48, 5).Value = '  +4.19%  '

$ws.Cells.Item(49, 4).Value = '''6.93'
$ws.Cells.Item(49, 5).Value = '  +6.57%  '

$ws.Cells.Item(50, 4).Value = '2.449.46'
$ws.Cells.Item(50, 5).Value = '  +11.65%  '

$ws.Cells.Item(51, 5).Value = '  +15.97%  '
